$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value
$updates = @(
    @{ Cell = 'D2'; Value = '29.417.46' },
    @{ Cell = 'D3'; Value = '1.848.36' },
    @{ Cell = 'E3'; Value = '  -0.28%  ' },
    @{ Cell = 'D4'; Value = '0.9990' },
    @{ Cell = 'E4'; Value = '  +0.01%  ' },
    @{ Cell = 'D5'; Value = '240.74' },
    @{ Cell = 'E5'; Value = '  -0.70%  ' },
    @{ Cell = 'D6'; Value = '0.6313' },
    @{ Cell = 'E7'; Value = '  +0.02%  ' },
    @{ Cell = 'E8'; Value = '  -0.16%  ' },
    @{ Cell = 'D9'; Value = '0.2957' },
    @{ Cell = 'E9'; Value = '  -1.02%  ' },
    @{ Cell = 'D10'; Value = '24.48' },
    @{ Cell = 'E10'; Value = '  -0.53%  ' },
    @{ Cell = 'D11'; Value = '0.07693' },
    @{ Cell = 'D12'; Value = '1.853.82' },
    @{ Cell = 'E12'; Value = '  -0.78%  ' },
    @{ Cell = 'E13'; Value = '  -0.60%  ' },
    @{ Cell = 'D14'; Value = '0.6852' },
    @{ Cell = 'E14'; Value = '  -1.01%  ' },
    @{ Cell = 'E15'; Value = '  +0.04%  ' },
    @{ Cell = 'D16'; Value = '83.12' },
    @{ Cell = 'E16'; Value = '  -0.59%  ' },
    @{ Cell = 'D17'; Value = '2.099.02' },
    @{ Cell = 'E17'; Value = '  -1.16%  ' },
    @{ Cell = 'D18'; Value = '6.129' },
    @{ Cell = 'E18'; Value = '  -2.21%  ' },
    @{ Cell = 'D19'; Value = '29.440.47' },
    @{ Cell = 'E19'; Value = '  -0.40%  ' },
    @{ Cell = 'D20'; Value = '228.05' },
    @{ Cell = 'E20'; Value = '  -2.47%  ' },
    @{ Cell = 'D21'; Value = '12.48' },
    @{ Cell = 'E21'; Value = '  -0.61%  ' },
    @{ Cell = 'D22'; Value = '0.9997' },
    @{ Cell = 'E22'; Value = '  -0.03%  ' },
    @{ Cell = 'D23'; Value = '7.551' },
    @{ Cell = 'E23'; Value = '  -1.50%  ' },
    @{ Cell = 'E24'; Value = '  +0.00%  ' },
    @{ Cell = 'D25'; Value = '157.02' },
    @{ Cell = 'E25'; Value = '  +1.05%  ' },
    @{ Cell = 'D26'; Value = '0.1393' },
    @{ Cell = 'E26'; Value = '  -0.52%  ' },
    @{ Cell = 'D27'; Value = '8.375' },
    @{ Cell = 'E27'; Value = '  -1.12%  ' },
    @{ Cell = 'D28'; Value = '17.68' },
    @{ Cell = 'E28'; Value = '  -0.25%  ' },
    @{ Cell = 'D29'; Value = '1.470' },
    @{ Cell = 'E29'; Value = '  -0.53%  ' },
    @{ Cell = 'D30'; Value = '1.264' },
    @{ Cell = 'E30'; Value = '  +0.80%  ' },
    @{ Cell = 'D31'; Value = '0.05711' },
    @{ Cell = 'E31'; Value = '  -2.22%  ' },
    @{ Cell = 'D32'; Value = '4.123' },
    @{ Cell = 'E32'; Value = '  -0.17%  ' },
    @{ Cell = 'D33'; Value = '4.025' },
    @{ Cell = 'E33'; Value = '  +0.02%  ' },
    @{ Cell = 'E34'; Value = '  -2.95%  ' },
    @{ Cell = 'D35'; Value = '1.156' },
    @{ Cell = 'E35'; Value = '  -1.26%  ' },
    @{ Cell = 'D36'; Value = '0.7129' },
    @{ Cell = 'E36'; Value = '  -1.18%  ' },
    @{ Cell = 'D37'; Value = '2.589' },
    @{ Cell = 'E37'; Value = '  +0.14%  ' },
    @{ Cell = 'D38'; Value = '1.249.92' },
    @{ Cell = 'E38'; Value = '  +0.65%  ' },
    @{ Cell = 'E39'; Value = '  +0.26%  ' },
    @{ Cell = 'D40'; Value = '2.779' },
    @{ Cell = 'E40'; Value = '  -0.73%  ' },
    @{ Cell = 'D41'; Value = '0.9081' },
    @{ Cell = 'E41'; Value = '  -0.05%  ' },
    @{ Cell = 'D42'; Value = '6.181' },
    @{ Cell = 'E42'; Value = '  +1.41%  ' },
    @{ Cell = 'D43'; Value = '1.000' },
    @{ Cell = 'E43'; Value = '  +0.08%  ' },
    @{ Cell = 'D44'; Value = '101.49' },
    @{ Cell = 'E44'; Value = '  -0.05%  ' },
    @{ Cell = 'D45'; Value = '66.17' },
    @{ Cell = 'E45'; Value = '  -2.39%  ' },
    @{ Cell = 'D46'; Value = '7.120' },
    @{ Cell = 'E46'; Value = '  -3.52%  ' },
    @{ Cell = 'B47'; Value = 'TheSandbox' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' },
    @{ Cell = 'D47'; Value = '0.4019' },
    @{ Cell = 'E47'; Value = '  -0.81%  ' },
    @{ Cell = 'B48'; Value = 'EnergySwap' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D48'; Value = '9.082' },
    @{ Cell = 'E48'; Value = '  -0.71%  ' },
    @{ Cell = 'B49'; Value = 'RenderToken' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D49'; Value = '1.680' },
    @{ Cell = 'E49'; Value = '  -1.65%  ' },
    @{ Cell = 'B50'; Value = 'Algorand' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D50'; Value = '0.1121' },
    @{ Cell = 'E50'; Value = '  +0.30%  ' },
    @{ Cell = 'B51'; Value = 'Cronos' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = 'D51'; Value = '0.05732' },
    @{ Cell = 'E51'; Value = '  -0.45%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
